$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D (Price) from Excel auto-converting numeric-looking
# text into actual numbers -- every value in D2:D51 is stored as text
# in the source workbook (t="inlineStr"), so force text format while
# we write the new values, then restore the original (default) style.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '36.698.82'
$ws.Range('E2').Value = '  +0.76%  '
$ws.Range('D3').Value = '1.967.84'
$ws.Range('E3').Value = '  +1.83%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '244.35'
$ws.Range('E5').Value = '  +1.13%  '
$ws.Range('D6').Value = '0.616'
$ws.Range('E6').Value = '  +1.26%  '
$ws.Range('D7').Value = '58.38'
$ws.Range('E7').Value = '  +2.87%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').Value = '0.373'
$ws.Range('E9').Value = '  +2.32%  '
$ws.Range('D10').Value = '0.0811'
$ws.Range('E10').Value = '  -1.36%  '
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('D12').Value = '22.12'
$ws.Range('E12').Value = '  +3.80%  '
$ws.Range('D13').Value = '2.256.16'
$ws.Range('E13').Value = '  +2.00%  '
$ws.Range('D14').Value = '0.823'
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('D15').Value = '13.71'
$ws.Range('E15').Value = '  +2.84%  '
$ws.Range('D16').Value = '5.27'
$ws.Range('E16').Value = '  +1.39%  '
$ws.Range('D17').Value = '1.966.30'
$ws.Range('E17').Value = '  -0.73%  '
$ws.Range('D18').Value = '36.713.45'
$ws.Range('E18').Value = '  +0.99%  '
$ws.Range('D19').Value = '69.66'
$ws.Range('E19').Value = '  +0.80%  '
$ws.Range('D20').Value = '0.0₃0861'
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('D21').Value = '5.10'
$ws.Range('E21').Value = '  +2.56%  '
$ws.Range('D22').Value = '228.12'
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').Value = '2.39'
$ws.Range('E24').Value = '  -1.91%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '2.34'
$ws.Range('E25').Value = '  +2.88%  '
$ws.Range('D26').Value = '9.33'
$ws.Range('E26').Value = '  +0.79%  '
$ws.Range('D27').Value = '160.94'
$ws.Range('E27').Value = '  -0.92%  '
$ws.Range('E28').Value = '  +10.79%  '
$ws.Range('D29').Value = '19.38'
$ws.Range('E29').Value = '  +0.77%  '
$ws.Range('D30').Value = '0.119'
$ws.Range('E30').Value = '  +1.47%  '
$ws.Range('D31').Value = '1.12'
$ws.Range('E31').Value = '  -0.99%  '
$ws.Range('D32').Value = '4.68'
$ws.Range('E32').Value = '  +1.02%  '
$ws.Range('D33').Value = '0.0620'
$ws.Range('E33').Value = '  -0.62%  '
$ws.Range('D34').Value = '4.25'
$ws.Range('E34').Value = '  -0.62%  '
$ws.Range('D35').Value = '6.23'
$ws.Range('E35').Value = '  +4.20%  '
$ws.Range('E36').Value = '  -0.23%  '
$ws.Range('D37').Value = '3.41'
$ws.Range('E37').Value = '  +17.32%  '
$ws.Range('D38').Value = '2.22'
$ws.Range('E38').Value = '  +4.09%  '
$ws.Range('E39').Value = '  -0.27%  '
$ws.Range('D40').Value = '0.0999'
$ws.Range('E40').Value = '  +3.22%  '
$ws.Range('D42').Value = '0.0212'
$ws.Range('E42').Value = '  +2.57%  '
$ws.Range('D43').Value = '1.17'
$ws.Range('E43').Value = '  -0.28%  '
$ws.Range('D44').Value = '16.05'
$ws.Range('E44').Value = '  +2.32%  '
$ws.Range('D45').Value = '1.05'
$ws.Range('E45').Value = '  +1.73%  '
$ws.Range('D46').Value = '1.358.78'
$ws.Range('E46').Value = '  +1.44%  '
$ws.Range('D47').Value = '87.50'
$ws.Range('E47').Value = '  +0.50%  '
$ws.Range('D48').Value = '7.15'
$ws.Range('E48').Value = '  -0.77%  '
$ws.Range('E49').Value = '  +1.39%  '
$ws.Range('D50').Value = '2.146.44'
$ws.Range('E50').Value = '  +1.88%  '
$ws.Range('D51').Value = '43.49'
$ws.Range('E51').Value = '  -4.63%  '

# Restore column D to its original (default/general) style now that the
# text values are safely written -- avoids leaving a text number format
# applied to the cells.
$ws.Range("D2:D51").Style = "Normal"

